# Auto-generated edit script: updates computed profit columns (H-N)
# for various crafting-leve rows across all eight job sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2318.5557
$ws.Range("I43").Value = 1190.5
$ws.Range("J43").Value = 2640.8572
$ws.Range("K43").Value = 1190.5
$ws.Range("L43").Value = 2640.8572
$ws.Range("M43").Value = -1121.5
$ws.Range("N43").Value = -2778.8572
$ws.Range("H80").Value = 529.02856
$ws.Range("I80").Value = 254.5
$ws.Range("J80").Value = 819.7059
$ws.Range("K80").Value = 763.5
$ws.Range("L80").Value = 2459.1177
$ws.Range("M80").Value = 234.5
$ws.Range("N80").Value = -4455.117700000001
$ws.Range("H83").Value = 529.02856
$ws.Range("I83").Value = 254.5
$ws.Range("J83").Value = 819.7059
$ws.Range("K83").Value = 2290.5
$ws.Range("L83").Value = 7377.3531
$ws.Range("M83").Value = 2701.5
$ws.Range("N83").Value = -17361.3531
$ws.Range("H133").Value = 51666.668
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 51666.668
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 51666.668
$ws.Range("N133").Value = -61786.668

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4771.039
$ws.Range("I32").Value = 4532.6904
$ws.Range("J32").Value = 5883.3335
$ws.Range("K32").Value = 4532.6904
$ws.Range("L32").Value = 5883.3335
$ws.Range("M32").Value = -4245.6904
$ws.Range("N32").Value = -6457.3335
$ws.Range("H63").Value = 13854494
$ws.Range("I63").Value = 17316396
$ws.Range("J63").Value = 6888.5
$ws.Range("K63").Value = 17316396
$ws.Range("L63").Value = 6888.5
$ws.Range("M63").Value = -17315710
$ws.Range("N63").Value = -8260.5
$ws.Range("H66").Value = 13854494
$ws.Range("I66").Value = 17316396
$ws.Range("J66").Value = 6888.5
$ws.Range("K66").Value = 86581980
$ws.Range("L66").Value = 34442.5
$ws.Range("M66").Value = -86578548
$ws.Range("N66").Value = -41306.5
$ws.Range("H97").Value = 871.8077
$ws.Range("I97").Value = 506.24
$ws.Range("J97").Value = 10011
$ws.Range("K97").Value = 506.24
$ws.Range("L97").Value = 10011
$ws.Range("M97").Value = -10.24000000000001
$ws.Range("N97").Value = -11003

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 36333.332
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 36333.332
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 36333.332
$ws.Range("N63").Value = -37705.332
$ws.Range("H66").Value = 36333.332
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 36333.332
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 108999.996
$ws.Range("N66").Value = -115863.996
$ws.Range("H69").Value = 32000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 32000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 32000
$ws.Range("N69").Value = -33622
$ws.Range("H72").Value = 32000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 32000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 96000
$ws.Range("N72").Value = -104112
$ws.Range("H86").Value = 2178.65
$ws.Range("I86").Value = 2000.2
$ws.Range("J86").Value = 2357.1
$ws.Range("K86").Value = 2000.2
$ws.Range("L86").Value = 2357.1
$ws.Range("M86").Value = -877.2
$ws.Range("N86").Value = -4603.1
$ws.Range("H87").Value = 44800
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 44800
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 44800
$ws.Range("N87").Value = -47296
$ws.Range("H89").Value = 2178.65
$ws.Range("I89").Value = 2000.2
$ws.Range("J89").Value = 2357.1
$ws.Range("K89").Value = 10001
$ws.Range("L89").Value = 11785.5
$ws.Range("M89").Value = -4385
$ws.Range("N89").Value = -23017.5
$ws.Range("H90").Value = 44800
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 44800
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 134400
$ws.Range("N90").Value = -146880
$ws.Range("H94").Value = 1312.7142
$ws.Range("I94").Value = 1241.8
$ws.Range("J94").Value = 1490
$ws.Range("K94").Value = 1241.8
$ws.Range("L94").Value = 1490
$ws.Range("M94").Value = -790.8
$ws.Range("N94").Value = -2392
$ws.Range("H95").Value = 39500
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 39500
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 39500
$ws.Range("N95").Value = -44992

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 18990
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 18990
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 18990
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -19874
$ws.Range("H56").Value = 40499.5
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 40499.5
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 40499.5
$ws.Range("N56").Value = -42189.5
$ws.Range("H57").Value = 46686.168
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 46686.168
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 45016
$ws.Range("M57").Value = 46686.168
$ws.Range("N57").Value = -47806.168
$ws.Range("H62").Value = 5133.3335
$ws.Range("I62").Value = 3950
$ws.Range("J62").Value = 7500
$ws.Range("K62").Value = 3950
$ws.Range("L62").Value = 7500
$ws.Range("M62").Value = -3326
$ws.Range("N62").Value = -8748
$ws.Range("H65").Value = 5133.3335
$ws.Range("I65").Value = 3950
$ws.Range("J65").Value = 7500
$ws.Range("K65").Value = 19750
$ws.Range("L65").Value = 37500
$ws.Range("M65").Value = -16630
$ws.Range("N65").Value = -43740
$ws.Range("H69").Value = 35335
$ws.Range("I69").Value = 13341
$ws.Range("J69").Value = 49997.668
$ws.Range("K69").Value = 13341
$ws.Range("L69").Value = 49997.668
$ws.Range("M69").Value = -12592
$ws.Range("N69").Value = -51495.668
$ws.Range("H72").Value = 35335
$ws.Range("I72").Value = 13341
$ws.Range("J72").Value = 49997.668
$ws.Range("K72").Value = 40023
$ws.Range("L72").Value = 149993.004
$ws.Range("M72").Value = -36279
$ws.Range("N72").Value = -157481.004
$ws.Range("H132").Value = 3192.0527
$ws.Range("I132").Value = 2218.6667
$ws.Range("J132").Value = 4860.7144
$ws.Range("K132").Value = 6656.000100000001
$ws.Range("L132").Value = 14582.1432
$ws.Range("M132").Value = -4126.000100000001
$ws.Range("N132").Value = -19642.1432

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 1330
$ws.Range("I47").Value = 995
$ws.Range("J47").Value = 2000
$ws.Range("K47").Value = 2985
$ws.Range("L47").Value = 6000
$ws.Range("M47").Value = -2554
$ws.Range("N47").Value = -6862
$ws.Range("H131").Value = 5495333
$ws.Range("I131").Value = 250000180
$ws.Range("J131").Value = 842.1573
$ws.Range("K131").Value = 750000540
$ws.Range("L131").Value = 2526.4719
$ws.Range("M131").Value = -749995500
$ws.Range("N131").Value = -12606.4719

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 62502530
$ws.Range("I80").Value = 83335704
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 83335704
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -83334706
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 62502530
$ws.Range("I83").Value = 83335704
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 416678520
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -416673528
$ws.Range("N83").Value = -24984

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2870
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 2783.3333
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 2783.3333
$ws.Range("M46").Value = -2812
$ws.Range("N46").Value = -3159.3333
$ws.Range("H132").Value = 4156
$ws.Range("I132").Value = 1504.4762
$ws.Range("J132").Value = 9218
$ws.Range("K132").Value = 4513.4286
$ws.Range("L132").Value = 27654
$ws.Range("M132").Value = -1983.4286
$ws.Range("N132").Value = -32714
$ws.Range("H133").Value = 33721.668
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 33721.668
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 33721.668
$ws.Range("N133").Value = -38781.668

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4510.5264
$ws.Range("I122").Value = 3108.5
$ws.Range("J122").Value = 6914
$ws.Range("K122").Value = 9325.5
$ws.Range("L122").Value = 20742
$ws.Range("M122").Value = -6875.5
$ws.Range("N122").Value = -25642
